$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 12025.25
$ws.Range("I86").Value = 12226.866
$ws.Range("K86").Value = 12226.866
$ws.Range("M86").Value = -11103.866
$ws.Range("H89").Value = 12025.25
$ws.Range("I89").Value = 12226.866
$ws.Range("K89").Value = 61134.33
$ws.Range("M89").Value = -55518.33
$ws.Range("H96").Value = 2939.5
$ws.Range("I96").Value = 2500
$ws.Range("J96").Value = 3379
$ws.Range("K96").Value = 7500
$ws.Range("L96").Value = 10137
$ws.Range("M96").Value = -6127
$ws.Range("N96").Value = -12883
$ws.Range("H100").Value = 16002
$ws.Range("I100").Value = 20000
$ws.Range("J100").Value = 14003
$ws.Range("K100").Value = 20000
$ws.Range("L100").Value = 14003
$ws.Range("M100").Value = -19459
$ws.Range("N100").Value = -15085
$ws.Range("H113").Value = 250000000
$ws.Range("I113").Value = 100000000
$ws.Range("J113").Value = 400000000
$ws.Range("K113").Value = 100000000
$ws.Range("L113").Value = 400000000
$ws.Range("M113").Value = -99996746
$ws.Range("N113").Value = -400006508
$ws.Range("H125").Value = 2709.5
$ws.Range("J125").Value = 4999
$ws.Range("L125").Value = 44991
$ws.Range("N125").Value = -49911
$ws.Range("H141").Value = 3242.5715
$ws.Range("I141").Value = 2938.6
$ws.Range("K141").Value = 8815.799999999999
$ws.Range("M141").Value = -3635.799999999999

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 80005
$ws.Range("J10").Value = 80005
$ws.Range("L10").Value = 80005
$ws.Range("N10").Value = -80345
$ws.Range("H40").Value = 23500
$ws.Range("I40").Value = 23500
$ws.Range("K40").Value = 23500
$ws.Range("M40").Value = -23324
$ws.Range("H61").Value = 20045764
$ws.Range("I61").Value = 27781896
$ws.Range("K61").Value = 27781896
$ws.Range("M61").Value = -27781684
$ws.Range("H74").Value = 13899893
$ws.Range("I74").Value = 20834884
$ws.Range("K74").Value = 20834884
$ws.Range("M74").Value = -20834010
$ws.Range("H77").Value = 13899893
$ws.Range("I77").Value = 20834884
$ws.Range("K77").Value = 104174420
$ws.Range("M77").Value = -104170052
$ws.Range("H110").Value = 1829.2
$ws.Range("I110").Value = 1726.8889
$ws.Range("J110").Value = 2750
$ws.Range("K110").Value = 1726.8889
$ws.Range("L110").Value = 2750
$ws.Range("M110").Value = 318.1111000000001
$ws.Range("N110").Value = -6840
$ws.Range("H132").Value = 4005.0435
$ws.Range("I132").Value = 3815.8
$ws.Range("K132").Value = 11447.4
$ws.Range("M132").Value = -8917.400000000001
$ws.Range("H136").Value = 20045764
$ws.Range("I136").Value = 27781896
$ws.Range("K136").Value = 83345688
$ws.Range("M136").Value = -83343138

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2857
$ws.Range("I105").Value = 1999.5
$ws.Range("J105").Value = 3200
$ws.Range("K105").Value = 1999.5
$ws.Range("L105").Value = 3200
$ws.Range("M105").Value = -252.5
$ws.Range("N105").Value = -6694
$ws.Range("H134").Value = 46390.086
$ws.Range("I134").Value = 1782.2858
$ws.Range("K134").Value = 5346.857400000001
$ws.Range("M134").Value = -2811.857400000001

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1115.625
$ws.Range("I58").Value = 819.3333
$ws.Range("J58").Value = 2004.5
$ws.Range("K58").Value = 819.3333
$ws.Range("L58").Value = 2004.5
$ws.Range("M58").Value = -616.3333
$ws.Range("N58").Value = -2410.5
$ws.Range("H86").Value = 100000
$ws.Range("I86").Value = 100000
$ws.Range("K86").Value = 100000
$ws.Range("M86").Value = -98877
$ws.Range("H89").Value = 100000
$ws.Range("I89").Value = 100000
$ws.Range("K89").Value = 500000
$ws.Range("M89").Value = -494384
$ws.Range("H99").Value = 2617.9443
$ws.Range("I99").Value = 2024.9231
$ws.Range("K99").Value = 2024.9231
$ws.Range("M99").Value = -526.9231
$ws.Range("H126").Value = 2617.9443
$ws.Range("I126").Value = 2024.9231
$ws.Range("K126").Value = 6074.7693
$ws.Range("M126").Value = -3604.7693
$ws.Range("H132").Value = 2456.318
$ws.Range("I132").Value = 2381.6316
$ws.Range("J132").Value = 2929.3333
$ws.Range("K132").Value = 7144.8948
$ws.Range("L132").Value = 8787.999899999999
$ws.Range("M132").Value = -4614.8948
$ws.Range("N132").Value = -13847.9999
$ws.Range("H134").Value = 629834.4399999999
$ws.Range("I134").Value = 1002312.6
$ws.Range("K134").Value = 3006937.8
$ws.Range("M134").Value = -3004402.8
$ws.Range("H136").Value = 1115.625
$ws.Range("I136").Value = 819.3333
$ws.Range("J136").Value = 2004.5
$ws.Range("K136").Value = 2457.9999
$ws.Range("L136").Value = 6013.5
$ws.Range("M136").Value = 92.0001000000002
$ws.Range("N136").Value = -11113.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 31716.355
$ws.Range("I129").Value = 672.3077
$ws.Range("J129").Value = 54137.055
$ws.Range("K129").Value = 2016.9231
$ws.Range("L129").Value = 162411.165
$ws.Range("M129").Value = 2983.0769
$ws.Range("N129").Value = -172411.165
$ws.Range("H131").Value = 6876.769
$ws.Range("I131").Value = 7731.25
$ws.Range("J131").Value = 5509.6
$ws.Range("K131").Value = 23193.75
$ws.Range("L131").Value = 16528.8
$ws.Range("M131").Value = -18153.75
$ws.Range("N131").Value = -26608.8

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 2528248.8
$ws.Range("J43").Value = 19990
$ws.Range("L43").Value = 19990
$ws.Range("N43").Value = -20376
$ws.Range("H100").Value = 1000
$ws.Range("I100").Value = 1000
$ws.Range("K100").Value = 1000
$ws.Range("M100").Value = -459
$ws.Range("H132").Value = 466271
$ws.Range("J132").Value = 106454.63
$ws.Range("L132").Value = 319363.89
$ws.Range("N132").Value = -324423.89
$ws.Range("H136").Value = 79785.734
$ws.Range("I136").Value = 2972.75
$ws.Range("K136").Value = 8918.25
$ws.Range("M136").Value = -6368.25

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H43").Value = 11927
$ws.Range("I43").Value = 11927
$ws.Range("K43").Value = 11927
$ws.Range("M43").Value = -11778
$ws.Range("H46").Value = 49500
$ws.Range("J46").Value = 49500
$ws.Range("L46").Value = 49500
$ws.Range("N46").Value = -49962
$ws.Range("H132").Value = 9239.147999999999
$ws.Range("I132").Value = 1933.5454
$ws.Range("K132").Value = 5800.6362
$ws.Range("M132").Value = -3270.6362
$ws.Range("H134").Value = 49500
$ws.Range("J134").Value = 49500
$ws.Range("L134").Value = 148500
$ws.Range("N134").Value = -153570
$ws.Range("H136").Value = 13772.75
$ws.Range("I136").Value = 1818.5454
$ws.Range("J136").Value = 28383.445
$ws.Range("K136").Value = 5455.6362
$ws.Range("L136").Value = 85150.33499999999
$ws.Range("M136").Value = -2905.6362
$ws.Range("N136").Value = -90250.33499999999
